$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsDeDe = $wb.Sheets.Item("de-de")

# "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" for the
# 56148062-90b8-49b6-a371-89468e1aa82c row (shared text across Overview!G
# and de-de!H): 2016-08-27 00:17:08 -> 2016-08-27 00:17:57
$wsOverview.Range("G3").Value = "2016-08-27 00:17:57"
$wsOverview.Range("G5").Value = "2016-08-27 00:17:57"
$wsDeDe.Range("H3").Value = "2016-08-27 00:17:57"
$wsDeDe.Range("H5").Value = "2016-08-27 00:17:57"

# Priority column: ht -> mt (shared between zh-cn and de-de rows)
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# zh-cn Correspond Handoff Datetime: 2016-08-27 00:16:59 -> 2016-08-27 00:17:53
$wsZhCn.Range("H3").Value = "2016-08-27 00:17:53"
$wsZhCn.Range("H5").Value = "2016-08-27 00:17:53"

# zh-cn Correspond Handback DateTime: 2016-08-27 00:17:28 -> 2016-08-27 00:18:14
$wsZhCn.Range("K3").Value = "2016-08-27 00:18:14"
$wsZhCn.Range("K5").Value = "2016-08-27 00:18:14"

# de-de Correspond Handback DateTime: 2016-08-27 00:17:34 -> 2016-08-27 00:18:22
$wsDeDe.Range("K3").Value = "2016-08-27 00:18:22"
$wsDeDe.Range("K5").Value = "2016-08-27 00:18:22"

$wb.Save()
